$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.636.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.111.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.92%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.33%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5260"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4507"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.85"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09013"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.104.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.48%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.824"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.025"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001175"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.94%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.36"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.290"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.710.29"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.83"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.80%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.353.93"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.53"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.528"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.633"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.41%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.014"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.891"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.20"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02643"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06824"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2308"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6867"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.283"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.97%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.320"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6418"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.80%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000358"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.246"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07282"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.49"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.13%  "
